$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text first so numeric-looking strings (e.g. "6.75") are not
# auto-converted to numbers by Excel when assigned below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.820.07"
$ws.Range("D3").Value = "2.917.23"
$ws.Range("E3").Value = "  -4.00%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "582.56"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "144.08"
$ws.Range("E6").Value = "  -5.91%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "2.916.36"
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "6.75"
$ws.Range("E10").Value = "  +6.77%  "
$ws.Range("D11").Value = "0.144"
$ws.Range("E11").Value = "  -4.11%  "
$ws.Range("E12").Value = "  -3.72%  "
$ws.Range("E13").Value = "  -3.49%  "
$ws.Range("D14").Value = "33.54"
$ws.Range("E14").Value = "  -5.11%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "3.401.05"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").Value = "60.783.15"
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("D18").Value = "6.73"
$ws.Range("E18").Value = "  -4.53%  "
$ws.Range("D19").Value = "2.918.03"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("D20").Value = "428.67"
$ws.Range("E20").Value = "  -4.95%  "
$ws.Range("D21").Value = "13.65"
$ws.Range("E21").Value = "  -4.34%  "
$ws.Range("D22").Value = "0.680"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").Value = "7.12"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "80.56"
$ws.Range("E24").Value = "  -3.29%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "11.00"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "2.21"
$ws.Range("E26").Value = "  -4.43%  "
$ws.Range("D27").Value = "11.84"
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "7.18"
$ws.Range("E30").Value = "  -4.46%  "
$ws.Range("D31").Value = "2.61"
$ws.Range("E31").Value = "  -3.30%  "
$ws.Range("D32").Value = "2.16"
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("D33").Value = "26.53"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("D35").Value = "0.0₃0877"
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("D37").Value = "5.64"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("D38").Value = "3.02"
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "49.80"
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.127"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "1.99"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("D42").Value = "8.67"
$ws.Range("E42").Value = "  -4.77%  "
$ws.Range("D43").Value = "0.296"
$ws.Range("E43").Value = "  -4.57%  "
$ws.Range("D44").Value = "41.43"
$ws.Range("E44").Value = "  -3.63%  "
$ws.Range("D45").Value = "378.61"
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("D46").Value = "0.0348"
$ws.Range("E46").Value = "  -3.19%  "
$ws.Range("D47").Value = "2.676.92"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("D48").Value = "132.06"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "24.51"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("E51").Value = "  -2.05%  "

# Restore the original (default) cell style for column D now that the text values
# are safely stored, so no visible formatting/style change is introduced.
$ws.Range("D2:D51").Style = "Normal"
